$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 34: num_customers 81 -> 83, retention_rate recalculated (83/2256)
$ws.Range("C34").Value = 83
$ws.Range("E34").Value = 0.03679078014184397

# Row 36: num_customers 138 -> 139, retention_rate recalculated (139/1930)
$ws.Range("C36").Value = 139
$ws.Range("E36").Value = 0.07202072538860103

# Row 37: num_customers 889 -> 892, cohort_size 889 -> 892 (retention_rate stays 1)
$ws.Range("C37").Value = 892
$ws.Range("D37").Value = 892
